$wb = $excel.ActiveWorkbook

$wsElectric = $wb.Worksheets.Item(1)
$wsGas = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Gas Network sheet: insert 3 new asset rows (Channelview, Elmax, Bacliff)
# right after the existing "Channel" row, which is simultaneously renamed to
# "Channel Energy". The rows that used to follow (Galveston, Bulldog, Free
# Port, Porth Arthur, Calcasieu Pass 2, Transmitter, Consumer) are pushed
# down by three rows.
# ---------------------------------------------------------------------------
$wsGas.Rows("15:17").Insert()

# Rename "Channel" -> "Channel Energy" (row 14, unaffected otherwise)
$wsGas.Cells.Item(14, 3).Value = "Channel Energy"

$nbsp = [char]0x00A0

# New row 15: Channelview
$wsGas.Cells.Item(15, 2).Value = "Producer "
$wsGas.Cells.Item(15, 3).Value = "Channelview"
$wsGas.Cells.Item(15, 4).Value = "Gas Power Station"
$wsGas.Cells.Item(15, 5).Value = 918
$wsGas.Cells.Item(15, 6).Value = "29.836952, -95.12174"
$wsGas.Cells.Item(15, 7).Value = 2002

# New row 16: Elmax
$wsGas.Cells.Item(16, 2).Value = "Producer "
$wsGas.Cells.Item(16, 3).Value = "Elmax"
$wsGas.Cells.Item(16, 4).Value = "Gas Power Station"
$wsGas.Cells.Item(16, 5).Value = 819
$wsGas.Cells.Item(16, 6).Value = "30.021922, -95.090427"
$wsGas.Cells.Item(16, 7).Value = 2028

# New row 17: Bacliff
$wsGas.Cells.Item(17, 2).Value = "Producer "
$wsGas.Cells.Item(17, 3).Value = "Bacliff"
$wsGas.Cells.Item(17, 4).Value = "Gas Power Station"
$wsGas.Cells.Item(17, 5).Value = 1036
$wsGas.Cells.Item(17, 6).Value = "29.49233, -94.98483" + $nbsp
$wsGas.Cells.Item(17, 7).Value = 2018

# ---------------------------------------------------------------------------
# View state: the workbook now opens on the "Gas Network" tab, with the
# "Electric Network" tab showing rows 15:21 selected and the "Gas Network"
# tab's selection resting on F21.
# ---------------------------------------------------------------------------
[void]$wsElectric.Activate()
[void]$wsElectric.Range("A15:XFD21").Select()

[void]$wsGas.Activate()
[void]$wsGas.Range("F21").Select()
